$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 0
$ws.Range("F6").Value = 15
$ws.Range("F14").Value = 8
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 0
$ws.Range("F23").Value = -1
$ws.Range("F25").Value = -3
$ws.Range("F26").Value = -2
$ws.Range("F28").Value = -4
$ws.Range("F32").Value = 0
